$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 170.93333
$ws.Range("I11").Value = 170.93333
$ws.Range("K11").Value = 170.93333
$ws.Range("M11").Value = -30.93333000000001
$ws.Range("H19").Value = 634.5
$ws.Range("I19").Value = 269
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 269
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -94
$ws.Range("N19").Value = -1350
$ws.Range("H32").Value = 6144.5
$ws.Range("J32").Value = 6144.5
$ws.Range("L32").Value = 6144.5
$ws.Range("N32").Value = -6796.5
$ws.Range("H33").Value = 11884.389
$ws.Range("I33").Value = 12571.647
$ws.Range("K33").Value = 12571.647
$ws.Range("M33").Value = -12342.647
$ws.Range("H40").Value = 6585.5713
$ws.Range("J40").Value = 4833.3335
$ws.Range("L40").Value = 4833.3335
$ws.Range("N40").Value = -5183.3335
$ws.Range("H62").Value = 4988.8887
$ws.Range("I62").Value = 825
$ws.Range("J62").Value = 6178.5713
$ws.Range("K62").Value = 825
$ws.Range("L62").Value = 6178.5713
$ws.Range("M62").Value = -201
$ws.Range("N62").Value = -7426.5713
$ws.Range("H65").Value = 4988.8887
$ws.Range("I65").Value = 825
$ws.Range("J65").Value = 6178.5713
$ws.Range("K65").Value = 4125
$ws.Range("L65").Value = 30892.8565
$ws.Range("M65").Value = -1005
$ws.Range("N65").Value = -37132.85649999999
$ws.Range("H98").Value = 3400.45
$ws.Range("I98").Value = 2597.5
$ws.Range("J98").Value = 6612.25
$ws.Range("K98").Value = 2597.5
$ws.Range("L98").Value = 6612.25
$ws.Range("M98").Value = -1099.5
$ws.Range("N98").Value = -9608.25
$ws.Range("H106").Value = 1799.25
$ws.Range("I106").Value = 1732.6666
$ws.Range("K106").Value = 1732.6666
$ws.Range("M106").Value = -1101.6666
$ws.Range("H111").Value = 4716.25
$ws.Range("I111").Value = 4516
$ws.Range("K111").Value = 13548
$ws.Range("M111").Value = -10481
$ws.Range("H116").Value = 7764.125
$ws.Range("I116").Value = 8000
$ws.Range("J116").Value = 7685.5
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 7685.5
$ws.Range("M116").Value = -4558
$ws.Range("N116").Value = -14569.5
$ws.Range("H118").Value = 587.8
$ws.Range("I118").Value = 587.8
$ws.Range("K118").Value = 1763.4
$ws.Range("M118").Value = -106.3999999999999
$ws.Range("H122").Value = 3400.45
$ws.Range("I122").Value = 2597.5
$ws.Range("J122").Value = 6612.25
$ws.Range("K122").Value = 7792.5
$ws.Range("L122").Value = 19836.75
$ws.Range("M122").Value = -5342.5
$ws.Range("N122").Value = -24736.75
$ws.Range("H131").Value = 4604.0527
$ws.Range("I131").Value = 3044.3845
$ws.Range("J131").Value = 7983.3335
$ws.Range("K131").Value = 9133.1535
$ws.Range("L131").Value = 23950.0005
$ws.Range("M131").Value = -4093.1535
$ws.Range("N131").Value = -34030.00049999999
$ws.Range("H132").Value = 3316.9429
$ws.Range("I132").Value = 2971.4688
$ws.Range("K132").Value = 8914.4064
$ws.Range("M132").Value = -6384.4064
$ws.Range("H137").Value = 4846.3335
$ws.Range("I137").Value = 1541.909
$ws.Range("K137").Value = 4625.727000000001
$ws.Range("M137").Value = -2075.727000000001
$ws.Range("H138").Value = 6862.3335
$ws.Range("I138").Value = 1309.3889
$ws.Range("J138").Value = 11027.042
$ws.Range("K138").Value = 3928.1667
$ws.Range("L138").Value = 33081.126
$ws.Range("M138").Value = 1211.8333
$ws.Range("N138").Value = -43361.126
$ws.Range("H141").Value = 8800.286
$ws.Range("I141").Value = 11012.5
$ws.Range("J141").Value = 3269.75
$ws.Range("K141").Value = 33037.5
$ws.Range("L141").Value = 9809.25
$ws.Range("M141").Value = -27857.5
$ws.Range("N141").Value = -20169.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4626.231
$ws.Range("I2").Value = 3704.125
$ws.Range("J2").Value = 6101.6
$ws.Range("K2").Value = 3704.125
$ws.Range("L2").Value = 6101.6
$ws.Range("M2").Value = -3591.125
$ws.Range("N2").Value = -6327.6
$ws.Range("H24").Value = 17500
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15748
$ws.Range("H32").Value = 6755
$ws.Range("I32").Value = 5964.5483
$ws.Range("J32").Value = 9477.666999999999
$ws.Range("K32").Value = 5964.5483
$ws.Range("L32").Value = 9477.666999999999
$ws.Range("M32").Value = -5677.5483
$ws.Range("N32").Value = -10051.667
$ws.Range("H61").Value = 23811600
$ws.Range("I61").Value = 29413670
$ws.Range("K61").Value = 29413670
$ws.Range("M61").Value = -29413458
$ws.Range("H63").Value = 9299.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 9299.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 9299.5
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -10671.5
$ws.Range("H66").Value = 9299.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 9299.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 46497.5
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -53361.5
$ws.Range("H74").Value = 20410740
$ws.Range("I74").Value = 52633588
$ws.Range("K74").Value = 52633588
$ws.Range("M74").Value = -52632714
$ws.Range("H77").Value = 20410740
$ws.Range("I77").Value = 52633588
$ws.Range("K77").Value = 263167940
$ws.Range("M77").Value = -263163572
$ws.Range("H100").Value = 17500
$ws.Range("J100").Value = 15000
$ws.Range("L100").Value = 15000
$ws.Range("N100").Value = -17164
$ws.Range("H110").Value = 4047.1072
$ws.Range("I110").Value = 3927.625
$ws.Range("J110").Value = 4206.4165
$ws.Range("K110").Value = 3927.625
$ws.Range("L110").Value = 4206.4165
$ws.Range("M110").Value = -1882.625
$ws.Range("N110").Value = -8296.416499999999
$ws.Range("H112").Value = 9346.75
$ws.Range("J112").Value = 9346.75
$ws.Range("L112").Value = 9346.75
$ws.Range("N112").Value = -12300.75
$ws.Range("H116").Value = 4626.231
$ws.Range("I116").Value = 3704.125
$ws.Range("J116").Value = 6101.6
$ws.Range("K116").Value = 3704.125
$ws.Range("L116").Value = 6101.6
$ws.Range("M116").Value = -1410.125
$ws.Range("N116").Value = -10689.6
$ws.Range("H124").Value = 34714.25
$ws.Range("J124").Value = 34714.25
$ws.Range("L124").Value = 34714.25
$ws.Range("N124").Value = -44534.25
$ws.Range("H130").Value = 46382.332
$ws.Range("J130").Value = 46382.332
$ws.Range("L130").Value = 46382.332
$ws.Range("N130").Value = -56422.332
$ws.Range("H132").Value = 22760788
$ws.Range("I132").Value = 1830.129
$ws.Range("J132").Value = 77032150
$ws.Range("K132").Value = 5490.387
$ws.Range("L132").Value = 231096450
$ws.Range("M132").Value = -2960.387
$ws.Range("N132").Value = -231101510
$ws.Range("H136").Value = 23811600
$ws.Range("I136").Value = 29413670
$ws.Range("K136").Value = 88241010
$ws.Range("M136").Value = -88238460

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4626.231
$ws.Range("I3").Value = 3704.125
$ws.Range("J3").Value = 6101.6
$ws.Range("K3").Value = 3704.125
$ws.Range("L3").Value = 6101.6
$ws.Range("M3").Value = -3590.125
$ws.Range("N3").Value = -6329.6
$ws.Range("H20").Value = 4429.857
$ws.Range("I20").Value = 3500
$ws.Range("K20").Value = 3500
$ws.Range("M20").Value = -3253
$ws.Range("H35").Value = 22487
$ws.Range("J35").Value = 39975
$ws.Range("L35").Value = 39975
$ws.Range("N35").Value = -40595
$ws.Range("H81").Value = 26352.8
$ws.Range("J81").Value = 26352.8
$ws.Range("L81").Value = 26352.8
$ws.Range("N81").Value = -28474.8
$ws.Range("H82").Value = 20033.273
$ws.Range("I82").Value = 4074
$ws.Range("J82").Value = 47962
$ws.Range("K82").Value = 4074
$ws.Range("L82").Value = 47962
$ws.Range("M82").Value = -3691
$ws.Range("N82").Value = -48728
$ws.Range("H84").Value = 26352.8
$ws.Range("J84").Value = 26352.8
$ws.Range("L84").Value = 79058.39999999999
$ws.Range("N84").Value = -89666.39999999999
$ws.Range("H85").Value = 20033.273
$ws.Range("I85").Value = 4074
$ws.Range("J85").Value = 47962
$ws.Range("K85").Value = 4074
$ws.Range("L85").Value = 47962
$ws.Range("M85").Value = -2748
$ws.Range("N85").Value = -50614
$ws.Range("H86").Value = 10125.25
$ws.Range("I86").Value = 11588.583
$ws.Range("J86").Value = 5735.25
$ws.Range("K86").Value = 11588.583
$ws.Range("L86").Value = 5735.25
$ws.Range("M86").Value = -10465.583
$ws.Range("N86").Value = -7981.25
$ws.Range("H89").Value = 10125.25
$ws.Range("I89").Value = 11588.583
$ws.Range("J89").Value = 5735.25
$ws.Range("K89").Value = 57942.915
$ws.Range("L89").Value = 28676.25
$ws.Range("M89").Value = -52326.915
$ws.Range("N89").Value = -39908.25
$ws.Range("H94").Value = 1801
$ws.Range("I94").Value = 1429.6428
$ws.Range("K94").Value = 1429.6428
$ws.Range("M94").Value = -978.6428000000001
$ws.Range("H99").Value = 1701.8096
$ws.Range("I99").Value = 1701.8096
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1701.8096
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -203.8096
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 6807.1035
$ws.Range("I105").Value = 17545.715
$ws.Range("K105").Value = 17545.715
$ws.Range("M105").Value = -15798.715
$ws.Range("H107").Value = 2595.5217
$ws.Range("I107").Value = 1735
$ws.Range("J107").Value = 8332.333000000001
$ws.Range("K107").Value = 1735
$ws.Range("L107").Value = 8332.333000000001
$ws.Range("M107").Value = 185
$ws.Range("N107").Value = -12172.333
$ws.Range("H110").Value = 76000
$ws.Range("J110").Value = 76000
$ws.Range("L110").Value = 76000
$ws.Range("N110").Value = -84180
$ws.Range("H130").Value = 58520.5
$ws.Range("I130").Value = 40000
$ws.Range("J130").Value = 67780.75
$ws.Range("K130").Value = 40000
$ws.Range("L130").Value = 67780.75
$ws.Range("M130").Value = -34980
$ws.Range("N130").Value = -77820.75
$ws.Range("H133").Value = 80420
$ws.Range("J133").Value = 80420
$ws.Range("L133").Value = 80420
$ws.Range("N133").Value = -90540
$ws.Range("H134").Value = 3140.04
$ws.Range("I134").Value = 3108.739
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 9326.217000000001
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -6791.217000000001
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 307.5
$ws.Range("I7").Value = 205
$ws.Range("J7").Value = 348.5
$ws.Range("K7").Value = 205
$ws.Range("L7").Value = 348.5
$ws.Range("M7").Value = -92
$ws.Range("N7").Value = -574.5
$ws.Range("H16").Value = 1494.1538
$ws.Range("I16").Value = 1382.8
$ws.Range("J16").Value = 1563.75
$ws.Range("K16").Value = 1382.8
$ws.Range("L16").Value = 1563.75
$ws.Range("M16").Value = -1095.8
$ws.Range("N16").Value = -2137.75
$ws.Range("H31").Value = 10874548
$ws.Range("I31").Value = 3155.1538
$ws.Range("J31").Value = 25007360
$ws.Range("K31").Value = 3155.1538
$ws.Range("L31").Value = 25007360
$ws.Range("M31").Value = -2860.1538
$ws.Range("N31").Value = -25007950
$ws.Range("H34").Value = 10874548
$ws.Range("I34").Value = 3155.1538
$ws.Range("J34").Value = 25007360
$ws.Range("K34").Value = 3155.1538
$ws.Range("L34").Value = 25007360
$ws.Range("M34").Value = -2953.1538
$ws.Range("N34").Value = -25007764
$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20884
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H94").Value = 1487.2
$ws.Range("I94").Value = 1632.6666
$ws.Range("J94").Value = 1461.5294
$ws.Range("K94").Value = 1632.6666
$ws.Range("L94").Value = 1461.5294
$ws.Range("M94").Value = -1181.6666
$ws.Range("N94").Value = -2363.5294
$ws.Range("H99").Value = 7080
$ws.Range("I99").Value = 7080
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7080
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -5582
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 18128.334
$ws.Range("I105").Value = 1504.5
$ws.Range("K105").Value = 1504.5
$ws.Range("M105").Value = 242.5
$ws.Range("H113").Value = 1494.1538
$ws.Range("I113").Value = 1382.8
$ws.Range("J113").Value = 1563.75
$ws.Range("K113").Value = 1382.8
$ws.Range("L113").Value = 1563.75
$ws.Range("M113").Value = 787.2
$ws.Range("N113").Value = -5903.75
$ws.Range("H122").Value = 2194780.5
$ws.Range("I122").Value = 1745.4762
$ws.Range("J122").Value = 17546026
$ws.Range("K122").Value = 5236.4286
$ws.Range("L122").Value = 52638078
$ws.Range("M122").Value = -2786.4286
$ws.Range("N122").Value = -52642978
$ws.Range("H126").Value = 7080
$ws.Range("I126").Value = 7080
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 21240
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -18770
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3804.5454
$ws.Range("I132").Value = 3785
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11355
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -8825
$ws.Range("N132").Value = -17060
$ws.Range("H134").Value = 3381.8572
$ws.Range("I134").Value = 2758.7273
$ws.Range("K134").Value = 8276.1819
$ws.Range("M134").Value = -5741.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17373962
$ws.Range("I4").Value = 38789964
$ws.Range("K4").Value = 116369892
$ws.Range("M4").Value = -116369780
$ws.Range("H11").Value = 136.66667
$ws.Range("I11").Value = 105
$ws.Range("K11").Value = 315
$ws.Range("M11").Value = -175
$ws.Range("H12").Value = 549.2273
$ws.Range("I12").Value = 533.8333
$ws.Range("J12").Value = 555
$ws.Range("K12").Value = 1601.4999
$ws.Range("L12").Value = 1665
$ws.Range("M12").Value = -1428.4999
$ws.Range("N12").Value = -2011
$ws.Range("H26").Value = 1692.7142
$ws.Range("J26").Value = 1819.8
$ws.Range("L26").Value = 5459.4
$ws.Range("N26").Value = -6035.4
$ws.Range("H56").Value = 8710.583000000001
$ws.Range("I56").Value = 8710.583000000001
$ws.Range("K56").Value = 8710.583000000001
$ws.Range("M56").Value = -8180.583000000001
$ws.Range("H59").Value = 4180
$ws.Range("I59").Value = 3633.3333
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 10899.9999
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = -10359.9999
$ws.Range("N59").Value = -16080
$ws.Range("H68").Value = 800
$ws.Range("J68").Value = 800
$ws.Range("L68").Value = 2400
$ws.Range("N68").Value = -4022
$ws.Range("H71").Value = 800
$ws.Range("J71").Value = 800
$ws.Range("L71").Value = 7200
$ws.Range("N71").Value = -15312
$ws.Range("H97").Value = 1538.2307
$ws.Range("J97").Value = 1374.3
$ws.Range("L97").Value = 4122.9
$ws.Range("N97").Value = -5114.9
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H107").Value = 959
$ws.Range("I107").Value = 251
$ws.Range("J107").Value = 1509.6666
$ws.Range("K107").Value = 753
$ws.Range("L107").Value = 4528.9998
$ws.Range("M107").Value = 1167
$ws.Range("N107").Value = -8368.9998
$ws.Range("H122").Value = 933.53845
$ws.Range("I122").Value = 354.5
$ws.Range("J122").Value = 1860
$ws.Range("K122").Value = 3190.5
$ws.Range("L122").Value = 16740
$ws.Range("M122").Value = -740.5
$ws.Range("N122").Value = -21640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1219801.8
$ws.Range("I2").Value = 4545505
$ws.Range("J2").Value = 377.2
$ws.Range("K2").Value = 4545505
$ws.Range("L2").Value = 377.2
$ws.Range("M2").Value = -4545392
$ws.Range("N2").Value = -603.2
$ws.Range("H36").Value = 5076.1113
$ws.Range("I36").Value = 3439
$ws.Range("K36").Value = 3439
$ws.Range("M36").Value = -2954
$ws.Range("H58").Value = 41800.4
$ws.Range("I58").Value = 38000
$ws.Range("J58").Value = 42750.5
$ws.Range("K58").Value = 38000
$ws.Range("L58").Value = 42750.5
$ws.Range("M58").Value = -37723
$ws.Range("N58").Value = -43304.5
$ws.Range("H70").Value = 7410.4443
$ws.Range("I70").Value = 6895
$ws.Range("J70").Value = 7474.875
$ws.Range("K70").Value = 6895
$ws.Range("L70").Value = 7474.875
$ws.Range("M70").Value = -6625
$ws.Range("N70").Value = -8014.875
$ws.Range("H73").Value = 7410.4443
$ws.Range("I73").Value = 6895
$ws.Range("J73").Value = 7474.875
$ws.Range("K73").Value = 6895
$ws.Range("L73").Value = 7474.875
$ws.Range("M73").Value = -5959
$ws.Range("N73").Value = -9346.875
$ws.Range("H97").Value = 963.3158
$ws.Range("J97").Value = 717.1429000000001
$ws.Range("L97").Value = 717.1429000000001
$ws.Range("N97").Value = -1709.1429
$ws.Range("H113").Value = 4796.2593
$ws.Range("I113").Value = 4055.4546
$ws.Range("K113").Value = 4055.4546
$ws.Range("M113").Value = -1885.4546
$ws.Range("H122").Value = 22729630
$ws.Range("I122").Value = 2094.2
$ws.Range("K122").Value = 6282.599999999999
$ws.Range("M122").Value = -3832.599999999999
$ws.Range("H126").Value = 73973.07000000001
$ws.Range("I126").Value = 118177.336
$ws.Range("K126").Value = 354532.008
$ws.Range("M126").Value = -352062.008
$ws.Range("H131").Value = 76333
$ws.Range("J131").Value = 76333
$ws.Range("L131").Value = 76333
$ws.Range("N131").Value = -86413
$ws.Range("H132").Value = 5093.8857
$ws.Range("I132").Value = 4880.355
$ws.Range("J132").Value = 6748.75
$ws.Range("K132").Value = 14641.065
$ws.Range("L132").Value = 20246.25
$ws.Range("M132").Value = -12111.065
$ws.Range("N132").Value = -25306.25
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5288.5557
$ws.Range("I7").Value = 5160.4287
$ws.Range("J7").Value = 5737
$ws.Range("K7").Value = 5160.4287
$ws.Range("L7").Value = 5737
$ws.Range("M7").Value = -5048.4287
$ws.Range("N7").Value = -5961
$ws.Range("H16").Value = 2110.4119
$ws.Range("I16").Value = 1086.9
$ws.Range("J16").Value = 3572.5715
$ws.Range("K16").Value = 1086.9
$ws.Range("L16").Value = 3572.5715
$ws.Range("M16").Value = -916.9000000000001
$ws.Range("N16").Value = -3912.5715
$ws.Range("H22").Value = 2342
$ws.Range("I22").Value = 2191.625
$ws.Range("J22").Value = 2642.75
$ws.Range("K22").Value = 2191.625
$ws.Range("L22").Value = 2642.75
$ws.Range("M22").Value = -1896.625
$ws.Range("N22").Value = -3232.75
$ws.Range("H27").Value = 2342
$ws.Range("I27").Value = 2191.625
$ws.Range("J27").Value = 2642.75
$ws.Range("K27").Value = 2191.625
$ws.Range("L27").Value = 2642.75
$ws.Range("M27").Value = -2084.625
$ws.Range("N27").Value = -2856.75
$ws.Range("H40").Value = 6149
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 32495
$ws.Range("J41").Value = 32495
$ws.Range("L41").Value = 32495
$ws.Range("N41").Value = -33371
$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H46").Value = 1621.8518
$ws.Range("I46").Value = 990
$ws.Range("K46").Value = 990
$ws.Range("M46").Value = -802
$ws.Range("H49").Value = 25000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H55").Value = 851.5
$ws.Range("I55").Value = 704.75
$ws.Range("J55").Value = 949.3333
$ws.Range("K55").Value = 704.75
$ws.Range("L55").Value = 949.3333
$ws.Range("M55").Value = -531.75
$ws.Range("N55").Value = -1295.3333
$ws.Range("H82").Value = 3299.2856
$ws.Range("I82").Value = 4097.3335
$ws.Range("K82").Value = 4097.3335
$ws.Range("M82").Value = -3736.3335
$ws.Range("H85").Value = 3299.2856
$ws.Range("I85").Value = 4097.3335
$ws.Range("K85").Value = 4097.3335
$ws.Range("M85").Value = -2849.3335
$ws.Range("H101").Value = 14281
$ws.Range("J101").Value = 14281
$ws.Range("L101").Value = 14281
$ws.Range("N101").Value = -20771
$ws.Range("H104").Value = 16059.75
$ws.Range("J104").Value = 16059.75
$ws.Range("L104").Value = 16059.75
$ws.Range("N104").Value = -23047.75
$ws.Range("H124").Value = 122499.5
$ws.Range("J124").Value = 122499.5
$ws.Range("L124").Value = 122499.5
$ws.Range("N124").Value = -132319.5
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H126").Value = 5288.5557
$ws.Range("I126").Value = 5160.4287
$ws.Range("J126").Value = 5737
$ws.Range("K126").Value = 15481.2861
$ws.Range("L126").Value = 17211
$ws.Range("M126").Value = -13011.2861
$ws.Range("N126").Value = -22151
$ws.Range("H127").Value = 80068.25
$ws.Range("J127").Value = 77270.86
$ws.Range("L127").Value = 77270.86
$ws.Range("N127").Value = -87190.86
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H132").Value = 3948.6365
$ws.Range("I132").Value = 2637.4285
$ws.Range("K132").Value = 7912.2855
$ws.Range("M132").Value = -5382.2855
$ws.Range("H133").Value = 116859.336
$ws.Range("J133").Value = 116859.336
$ws.Range("L133").Value = 116859.336
$ws.Range("N133").Value = -121919.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H81").Value = 6806
$ws.Range("I81").Value = 6499.5
$ws.Range("J81").Value = 6867.3
$ws.Range("K81").Value = 12999
$ws.Range("L81").Value = 13734.6
$ws.Range("M81").Value = -11938
$ws.Range("N81").Value = -15856.6
$ws.Range("H84").Value = 6806
$ws.Range("I84").Value = 6499.5
$ws.Range("J84").Value = 6867.3
$ws.Range("K84").Value = 64995
$ws.Range("L84").Value = 68673
$ws.Range("M84").Value = -59691
$ws.Range("N84").Value = -79281
$ws.Range("H103").Value = 18405.75
$ws.Range("J103").Value = 18405.75
$ws.Range("L103").Value = 18405.75
$ws.Range("N103").Value = -20749.75
$ws.Range("H124").Value = 5037500
$ws.Range("J124").Value = 5037500
$ws.Range("L124").Value = 5037500
$ws.Range("N124").Value = -5047320
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H133").Value = 72858.8
$ws.Range("J133").Value = 72858.8
$ws.Range("L133").Value = 72858.8
$ws.Range("N133").Value = -82978.8
$ws.Range("H135").Value = 25035250
$ws.Range("J135").Value = 25035250
$ws.Range("L135").Value = 25035250
$ws.Range("N135").Value = -25045390
$ws.Range("H137").Value = 99998.336
$ws.Range("J137").Value = 99998.336
$ws.Range("L137").Value = 99998.336
$ws.Range("N137").Value = -110198.336
